# Applies the committed changes:
#  - fix two data-entry typos in the "Secteur" column (shared strings)
#  - widen column B ("Secteur") so the longer label is readable
#  - update the last active selection to the cells that were corrected
#  - shrink the saved window height

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the two mis-typed "Secteur" values -----------------------
# "Torra di Murtella" -> "A Torra di Murtella" (rows 18-21)
$ws.Range("B18:B21").Value = "A Torra di Murtella"
# "Maffalcu" -> "Malfalcu" (rows 30-33)
$ws.Range("B30:B33").Value = "Malfalcu"

# --- Widen column B (Secteur) to fit the corrected, longer label ------
$ws.Columns.Item(2).ColumnWidth = 38.6190476190476

# --- Move/extend the selection to the cells that were just corrected --
$ws.Range("B18:B21").Select() | Out-Null

# --- Shrink the recorded window height --------------------------------
$win = $excel.ActiveWindow
$win.Height = 12180
$win.Width = 27945
